$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")

# Rows 7-18: "Mass" quantities (previously interleaved with "Weight" rows)
$ws.Range("A7").Value = "Take-Off Mass"
$ws.Range("B7").Value = "kg"
$ws.Range("C7").Value = 57925.350866929904

$ws.Range("A8").Value = "Maximum Landing Mass"
$ws.Range("B8").Value = "kg"
$ws.Range("C8").Value = 52132.81578023691

$ws.Range("A9").Value = "Maximum Passengers Mass"
$ws.Range("B9").Value = "kg"
$ws.Range("C9").Value = 12870.0

$ws.Range("A10").Value = "Fuel Mass"
$ws.Range("B10").Value = "kg"
$ws.Range("C10").Value = 11959.066679280662

$ws.Range("A11").Value = "Crew Mass"
$ws.Range("B11").Value = "kg"
$ws.Range("C11").Value = 459.08729100000005

$ws.Range("A12").Value = "Maximum Zero Fuel Mass"
$ws.Range("B12").Value = "kg"
$ws.Range("C12").Value = 45966.284187649246

$ws.Range("A13").Value = "Zero Fuel Mass"
$ws.Range("B13").Value = "kg"
$ws.Range("C13").Value = 45966.284187649246

$ws.Range("A14").Value = "Operating Empty Mass"
$ws.Range("B14").Value = "kg"
$ws.Range("C14").Value = 33096.284187649246

$ws.Range("A15").Value = "Empty Mass"
$ws.Range("B15").Value = "kg"
$ws.Range("C15").Value = 32367.196896649246

$ws.Range("A16").Value = "Manufacturer Empty Mass"
$ws.Range("B16").Value = "kg"
$ws.Range("C16").Value = 31516.98689664925

$ws.Range("A17").Value = "Operating Item Mass"
$ws.Range("B17").Value = "kg"
$ws.Range("C17").Value = 1120.21

$ws.Range("A18").Value = "Trapped Fuel Oil Mass"
$ws.Range("B18").Value = "kg"
$ws.Range("C18").Value = 270.0

# Row 19: blank separator row (single space in column A, no B/C)
$ws.Range("A19").Value = " "
$ws.Range("B19").ClearContents()
$ws.Range("C19").ClearContents()

# Rows 20-32: "Weight" quantities
$ws.Range("A20").Value = "Maximum Take-Off Weight"
$ws.Range("B20").Value = "N"
$ws.Range("C20").Value = 562518.9322633445

$ws.Range("A21").Value = "Take-Off Weight"
$ws.Range("B21").Value = "N"
$ws.Range("C21").Value = 568053.642079178

$ws.Range("A22").Value = "Maximum Landing Weight"
$ws.Range("B22").Value = "N"
$ws.Range("C22").Value = 511248.2778712602

$ws.Range("A23").Value = "Maximum Passengers Weight"
$ws.Range("B23").Value = "N"
$ws.Range("C23").Value = 126211.58549999996

$ws.Range("A24").Value = "Fuel Weight"
$ws.Range("B24").Value = "N"
$ws.Range("C24").Value = 117278.38125036767

$ws.Range("A25").Value = "Crew Weight"
$ws.Range("B25").Value = "N"
$ws.Range("C25").Value = 4502.108382285149

$ws.Range("A26").Value = "Maximum Zero Fuel Weight"
$ws.Range("B26").Value = "N"
$ws.Range("C26").Value = 450775.2608288104

$ws.Range("A27").Value = "Zero Fuel Weight"
$ws.Range("B27").Value = "N"
$ws.Range("C27").Value = 450775.2608288104

$ws.Range("A28").Value = "Operating Empty Weight"
$ws.Range("B28").Value = "N"
$ws.Range("C28").Value = 324563.67532881047

$ws.Range("A29").Value = "Empty Weight"
$ws.Range("B29").Value = "N"
$ws.Range("C29").Value = 317413.77144652524

$ws.Range("A30").Value = "Manufacturer Empty Weight"
$ws.Range("B30").Value = "N"
$ws.Range("C30").Value = 309076.0595500253

$ws.Range("A31").Value = "Operating Item Weight"
$ws.Range("B31").Value = "N"
$ws.Range("C31").Value = 10985.507396499997

$ws.Range("A32").Value = "Trapped Fuel Oil Weight"
$ws.Range("B32").Value = "N"
$ws.Range("C32").Value = 2647.7954999999993

# Row 33 no longer exists (the duplicate "Operating Empty Weight" row is removed)
$ws.Range("A33:C33").ClearContents()
